$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data was added: insert a new record as row 5,
# pushing the existing rows 5-39 down to 6-40.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 44532
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100103
$ws.Range("H5").Value = "Frutos de hueso (carozo)"
$ws.Range("I5").Value = 100103001
$ws.Range("J5").Value = "Cereza"
$ws.Range("K5").Value = "Santina"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20500
$ws.Range("Q5").Value = "$/caja 15 kilos"
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 1367
$ws.Range("T5").Value = 15
